# Edit script implementing the diff changes for Введение.docx
$d = $word.ActiveDocument

function Replace-Text {
    param(
        [string]$Find,
        [string]$Replace,
        [bool]$MatchWildcards = $false
    )
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Replacement.ClearFormatting()
    $ok = $range.Find.Execute($Find, $true, $false, $MatchWildcards, $false, $false, $true, 1, $false, $Replace, 2)
    if (-not $ok) {
        Write-Host "WARNING: replace failed for: $Find"
    }
    return $ok
}

# 1. Split title paragraph into two paragraphs
Replace-Text "к дипломной работе по теме «Программное средство для расчета химического состава образцов»" "к дипломной работе по теме:^p«Программное средство для расчета химического состава образцов»"

# Fix run formatting (szCs=28) on the new second paragraph's run, matching sibling paragraphs
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Программное средство для расчета химического состава образцов*") {
        $p.Range.Font.SizeBi = 14
        break
    }
}

# 2. "...калибровочных данных, различных настроек и пр." -> append ";"
Replace-Text "различных настроек и пр." "различных настроек и пр.;"

# 3. "...выполненных по заданным схемах, в заданный интервал времени." -> end with ";"
Replace-Text "выполненных по заданным схемах, в заданный интервал времени." "выполненных по заданным схемах, в заданный интервал времени;"

# 4. Insert "программные " before "средства:"
Replace-Text "следующие средства:" "следующие программные средства:"

# 5. Append ";" after "15.7.5"
Replace-Text "Community 2017 версия 15.7.5" "Community 2017 версия 15.7.5;"

# 6. WPF Framework trailing space -> ";"
Replace-Text "WPF Framework " "WPF Framework;"

# 7. " в качестве средства обеспечения доступа к БД" -> "...к базе данных;"
Replace-Text "в качестве средства обеспечения доступа к БД" "в качестве средства обеспечения доступа к базе данных;"

# 8. Big rewrite of the LocalDB justification sentence (part 1: up to "инфраструктура")
Replace-Text "так как при подключении с использованием строки подключения автоматически создается и запускается требуемая инфраструктура" "наиболее приемлемый вариант, принимая в расчет тот факт, что в геохимических лабораториях, как правило, отсутствует штатный сотрудник, который следил бы за состоянием и работой полноценного MS SQL Server'а, а также особенность данной версии, состоящей в автоматическом создании и запуске требуемой инфраструктуры"

# 9. Insert clause after "конфигурационных задач"
Replace-Text "сложных конфигурационных задач " "сложных конфигурационных задач, в случае если подключение осуществляется через строку подключения "

# 10. Final run after citation sdt: " наиболее приемлемый вариант " -> "."
Replace-Text " наиболее приемлемый вариант " "."
